$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" column, reusing the header formatting from the
# existing last header cell (G1) so H1 picks up the same style (bold,
# centered, bordered) as the other header cells.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 1
